$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix D7 value (was TRUE, should now be FALSE)
$ws.Range("D7").Value = $false

# Add new row 9 with test case data
$ws.Range("A9").Value = "TEST_CASE_8"
$ws.Range("B9").Value = "Verifica upload"
$ws.Range("C9").Value = "Apri il sito https://www.zamzar.com/ e clicca sull'icona ""triangolo rovesciato"" di fianco a""Choose file"" e seleziona da URL, incolla https://avatars.githubusercontent.com/u/192012301?s=48&v=4 e verifica l'upload"
$ws.Range("D9").Value = $true
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "web"
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
